$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 2).Value = 0.298
$ws.Cells.Item(4, 5).Value = 0.176
$ws.Cells.Item(4, 6).Value = 0.028
$ws.Cells.Item(4, 7).Value = 0.167
$ws.Cells.Item(4, 8).Value = 0.209
$ws.Cells.Item(4, 9).Value = 0.025
$ws.Cells.Item(4, 10).Value = 0.159
$ws.Cells.Item(4, 11).Value = 0.322
$ws.Cells.Item(4, 12).Value = 0.099
$ws.Cells.Item(4, 13).Value = 0.314
$ws.Cells.Item(4, 14).Value = 0.26
$ws.Cells.Item(4, 16).Value = 0.147
$ws.Cells.Item(4, 17).Value = 0.5
$ws.Cells.Item(4, 18).Value = 0.221
$ws.Cells.Item(4, 19).Value = 0.47
$ws.Cells.Item(4, 20).Value = 0.256
$ws.Cells.Item(4, 22).Value = 0.291
$ws.Cells.Item(4, 23).Value = 0.253
$ws.Cells.Item(4, 24).Value = 0.042
$ws.Cells.Item(4, 25).Value = 0.204
$ws.Cells.Item(4, 26).Value = 0.451
$ws.Cells.Item(4, 27).Value = 0.127
$ws.Cells.Item(4, 28).Value = 0.357
$ws.Cells.Item(4, 29).Value = 0.118
$ws.Cells.Item(4, 31).Value = 0.081
$ws.Cells.Item(4, 32).Value = 0.717
$ws.Cells.Item(4, 34).Value = 0.322
$ws.Cells.Item(4, 35).Value = 0.674
$ws.Cells.Item(4, 36).Value = 0.161
$ws.Cells.Item(4, 37).Value = 0.402
$ws.Cells.Item(4, 38).Value = 0.695
$ws.Cells.Item(4, 41).Value = 0.695

# Row 5
$ws.Cells.Item(5, 2).Value = 0.837
$ws.Cells.Item(5, 3).Value = 0.136
$ws.Cells.Item(5, 4).Value = 0.369
$ws.Cells.Item(5, 5).Value = 0.698
$ws.Cells.Item(5, 6).Value = 0.211
$ws.Cells.Item(5, 7).Value = 0.459
$ws.Cells.Item(5, 8).Value = 0.837
$ws.Cells.Item(5, 9).Value = 0.136
$ws.Cells.Item(5, 10).Value = 0.369
$ws.Cells.Item(5, 11).Value = 0.628
$ws.Cells.Item(5, 12).Value = 0.234
$ws.Cells.Item(5, 13).Value = 0.483
$ws.Cells.Item(5, 14).Value = 0.814
$ws.Cells.Item(5, 15).Value = 0.151
$ws.Cells.Item(5, 16).Value = 0.389
$ws.Cells.Item(5, 17).Value = 0.558
$ws.Cells.Item(5, 18).Value = 0.247
$ws.Cells.Item(5, 19).Value = 0.497
$ws.Cells.Item(5, 20).Value = 0.535
$ws.Cells.Item(5, 21).Value = 0.249
$ws.Cells.Item(5, 22).Value = 0.499
$ws.Cells.Item(5, 23).Value = 0.767
$ws.Cells.Item(5, 24).Value = 0.178
$ws.Cells.Item(5, 25).Value = 0.422
$ws.Cells.Item(5, 26).Value = 0.837
$ws.Cells.Item(5, 27).Value = 0.136
$ws.Cells.Item(5, 28).Value = 0.369
$ws.Cells.Item(5, 29).Value = 0.721
$ws.Cells.Item(5, 30).Value = 0.201
$ws.Cells.Item(5, 31).Value = 0.449
$ws.Cells.Item(5, 32).Value = 0.953
$ws.Cells.Item(5, 33).Value = 0.044
$ws.Cells.Item(5, 34).Value = 0.211
$ws.Cells.Item(5, 35).Value = 0.791
$ws.Cells.Item(5, 36).Value = 0.165
$ws.Cells.Item(5, 37).Value = 0.407
$ws.Cells.Item(5, 38).Value = 0.93
$ws.Cells.Item(5, 39).Value = 0.065
$ws.Cells.Item(5, 40).Value = 0.255
$ws.Cells.Item(5, 41).Value = 0.891

# Row 6
$ws.Cells.Item(6, 2).Value = 0.44
$ws.Cells.Item(6, 5).Value = 0.281
$ws.Cells.Item(6, 8).Value = 0.334
$ws.Cells.Item(6, 11).Value = 0.426
$ws.Cells.Item(6, 14).Value = 0.394
$ws.Cells.Item(6, 17).Value = 0.527
$ws.Cells.Item(6, 20).Value = 0.346
$ws.Cells.Item(6, 23).Value = 0.38
$ws.Cells.Item(6, 26).Value = 0.586
$ws.Cells.Item(6, 29).Value = 0.203
$ws.Cells.Item(6, 32).Value = 0.818
$ws.Cells.Item(6, 35).Value = 0.728
$ws.Cells.Item(6, 38).Value = 0.796
$ws.Cells.Item(6, 41).Value = 0.781

# Row 7
$ws.Cells.Item(7, 2).Value = 0.615
$ws.Cells.Item(7, 5).Value = 0.438
$ws.Cells.Item(7, 8).Value = 0.523
$ws.Cells.Item(7, 11).Value = 0.528
$ws.Cells.Item(7, 14).Value = 0.571
$ws.Cells.Item(7, 17).Value = 0.545
$ws.Cells.Item(7, 20).Value = 0.439
$ws.Cells.Item(7, 23).Value = 0.545
$ws.Cells.Item(7, 26).Value = 0.715
$ws.Cells.Item(7, 29).Value = 0.357
$ws.Cells.Item(7, 32).Value = 0.894
$ws.Cells.Item(7, 35).Value = 0.764
$ws.Cells.Item(7, 38).Value = 0.871
$ws.Cells.Item(7, 41).Value = 0.843

# Row 8
$ws.Cells.Item(8, 2).Value = 0.753
$ws.Cells.Item(8, 3).Value = 0.141
$ws.Cells.Item(8, 4).Value = 0.376
$ws.Cells.Item(8, 5).Value = 0.589
$ws.Cells.Item(8, 7).Value = 0.431
$ws.Cells.Item(8, 8).Value = 0.729
$ws.Cells.Item(8, 9).Value = 0.142
$ws.Cells.Item(8, 10).Value = 0.377
$ws.Cells.Item(8, 11).Value = 0.556
$ws.Cells.Item(8, 12).Value = 0.207
$ws.Cells.Item(8, 13).Value = 0.455
$ws.Cells.Item(8, 14).Value = 0.735
$ws.Cells.Item(8, 15).Value = 0.151
$ws.Cells.Item(8, 16).Value = 0.388
$ws.Cells.Item(8, 17).Value = 0.532
$ws.Cells.Item(8, 18).Value = 0.233
$ws.Cells.Item(8, 19).Value = 0.482
$ws.Cells.Item(8, 20).Value = 0.466
$ws.Cells.Item(8, 21).Value = 0.21
$ws.Cells.Item(8, 23).Value = 0.693
$ws.Cells.Item(8, 24).Value = 0.169
$ws.Cells.Item(8, 25).Value = 0.411
$ws.Cells.Item(8, 26).Value = 0.771
$ws.Cells.Item(8, 27).Value = 0.138
$ws.Cells.Item(8, 28).Value = 0.371
$ws.Cells.Item(8, 29).Value = 0.6
$ws.Cells.Item(8, 31).Value = 0.431
$ws.Cells.Item(8, 32).Value = 0.882
$ws.Cells.Item(8, 33).Value = 0.061
$ws.Cells.Item(8, 34).Value = 0.246
$ws.Cells.Item(8, 35).Value = 0.782
$ws.Cells.Item(8, 36).Value = 0.165
$ws.Cells.Item(8, 37).Value = 0.406
$ws.Cells.Item(8, 38).Value = 0.896
$ws.Cells.Item(8, 39).Value = 0.072
$ws.Cells.Item(8, 40).Value = 0.268
$ws.Cells.Item(8, 41).Value = 0.853

# Row 9
$ws.Cells.Item(9, 2).Value = 0.651
$ws.Cells.Item(9, 3).Value = 0.227
$ws.Cells.Item(9, 4).Value = 0.477
$ws.Cells.Item(9, 5).Value = 0.465
$ws.Cells.Item(9, 6).Value = 0.249
$ws.Cells.Item(9, 7).Value = 0.499
$ws.Cells.Item(9, 8).Value = 0.605
$ws.Cells.Item(9, 9).Value = 0.239
$ws.Cells.Item(9, 10).Value = 0.489
$ws.Cells.Item(9, 11).Value = 0.465
$ws.Cells.Item(9, 12).Value = 0.249
$ws.Cells.Item(9, 13).Value = 0.499
$ws.Cells.Item(9, 14).Value = 0.628
$ws.Cells.Item(9, 15).Value = 0.234
$ws.Cells.Item(9, 16).Value = 0.483
$ws.Cells.Item(9, 17).Value = 0.488
$ws.Cells.Item(9, 20).Value = 0.372
$ws.Cells.Item(9, 21).Value = 0.234
$ws.Cells.Item(9, 22).Value = 0.483
$ws.Cells.Item(9, 23).Value = 0.581
$ws.Cells.Item(9, 24).Value = 0.243
$ws.Cells.Item(9, 25).Value = 0.493
$ws.Cells.Item(9, 26).Value = 0.674
$ws.Cells.Item(9, 27).Value = 0.22
$ws.Cells.Item(9, 28).Value = 0.469
$ws.Cells.Item(9, 29).Value = 0.488
$ws.Cells.Item(9, 30).Value = 0.25
$ws.Cells.Item(9, 31).Value = 0.5
$ws.Cells.Item(9, 32).Value = 0.767
$ws.Cells.Item(9, 33).Value = 0.178
$ws.Cells.Item(9, 34).Value = 0.422
$ws.Cells.Item(9, 35).Value = 0.767
$ws.Cells.Item(9, 36).Value = 0.178
$ws.Cells.Item(9, 37).Value = 0.422
$ws.Cells.Item(9, 38).Value = 0.837
$ws.Cells.Item(9, 39).Value = 0.136
$ws.Cells.Item(9, 40).Value = 0.369
$ws.Cells.Item(9, 41).Value = 0.79

# Row 10
$ws.Cells.Item(10, 2).Value = 0.791
$ws.Cells.Item(10, 3).Value = 0.165
$ws.Cells.Item(10, 4).Value = 0.407
$ws.Cells.Item(10, 5).Value = 0.628
$ws.Cells.Item(10, 6).Value = 0.234
$ws.Cells.Item(10, 7).Value = 0.483
$ws.Cells.Item(10, 8).Value = 0.767
$ws.Cells.Item(10, 9).Value = 0.178
$ws.Cells.Item(10, 10).Value = 0.422
$ws.Cells.Item(10, 11).Value = 0.628
$ws.Cells.Item(10, 12).Value = 0.234
$ws.Cells.Item(10, 13).Value = 0.483
$ws.Cells.Item(10, 14).Value = 0.791
$ws.Cells.Item(10, 15).Value = 0.165
$ws.Cells.Item(10, 16).Value = 0.407
$ws.Cells.Item(10, 17).Value = 0.558
$ws.Cells.Item(10, 18).Value = 0.247
$ws.Cells.Item(10, 19).Value = 0.497
$ws.Cells.Item(10, 20).Value = 0.535
$ws.Cells.Item(10, 21).Value = 0.249
$ws.Cells.Item(10, 22).Value = 0.499
$ws.Cells.Item(10, 23).Value = 0.767
$ws.Cells.Item(10, 24).Value = 0.178
$ws.Cells.Item(10, 25).Value = 0.422
$ws.Cells.Item(10, 26).Value = 0.837
$ws.Cells.Item(10, 27).Value = 0.136
$ws.Cells.Item(10, 28).Value = 0.369
$ws.Cells.Item(10, 29).Value = 0.605
$ws.Cells.Item(10, 30).Value = 0.239
$ws.Cells.Item(10, 31).Value = 0.489
$ws.Cells.Item(10, 32).Value = 0.953
$ws.Cells.Item(10, 33).Value = 0.044
$ws.Cells.Item(10, 34).Value = 0.211
$ws.Cells.Item(10, 35).Value = 0.791
$ws.Cells.Item(10, 36).Value = 0.165
$ws.Cells.Item(10, 37).Value = 0.407
$ws.Cells.Item(10, 38).Value = 0.93
$ws.Cells.Item(10, 39).Value = 0.065
$ws.Cells.Item(10, 40).Value = 0.255
$ws.Cells.Item(10, 41).Value = 0.891

# Row 11
$ws.Cells.Item(11, 2).Value = 0.837
$ws.Cells.Item(11, 3).Value = 0.136
$ws.Cells.Item(11, 4).Value = 0.369
$ws.Cells.Item(11, 5).Value = 0.698
$ws.Cells.Item(11, 6).Value = 0.211
$ws.Cells.Item(11, 7).Value = 0.459
$ws.Cells.Item(11, 8).Value = 0.837
$ws.Cells.Item(11, 9).Value = 0.136
$ws.Cells.Item(11, 10).Value = 0.369
$ws.Cells.Item(11, 11).Value = 0.628
$ws.Cells.Item(11, 12).Value = 0.234
$ws.Cells.Item(11, 13).Value = 0.483
$ws.Cells.Item(11, 14).Value = 0.814
$ws.Cells.Item(11, 15).Value = 0.151
$ws.Cells.Item(11, 16).Value = 0.389
$ws.Cells.Item(11, 17).Value = 0.558
$ws.Cells.Item(11, 18).Value = 0.247
$ws.Cells.Item(11, 19).Value = 0.497
$ws.Cells.Item(11, 20).Value = 0.535
$ws.Cells.Item(11, 21).Value = 0.249
$ws.Cells.Item(11, 22).Value = 0.499
$ws.Cells.Item(11, 23).Value = 0.767
$ws.Cells.Item(11, 24).Value = 0.178
$ws.Cells.Item(11, 25).Value = 0.422
$ws.Cells.Item(11, 26).Value = 0.837
$ws.Cells.Item(11, 27).Value = 0.136
$ws.Cells.Item(11, 28).Value = 0.369
$ws.Cells.Item(11, 29).Value = 0.651
$ws.Cells.Item(11, 30).Value = 0.227
$ws.Cells.Item(11, 31).Value = 0.477
$ws.Cells.Item(11, 32).Value = 0.953
$ws.Cells.Item(11, 33).Value = 0.044
$ws.Cells.Item(11, 34).Value = 0.211
$ws.Cells.Item(11, 35).Value = 0.791
$ws.Cells.Item(11, 36).Value = 0.165
$ws.Cells.Item(11, 37).Value = 0.407
$ws.Cells.Item(11, 38).Value = 0.93
$ws.Cells.Item(11, 39).Value = 0.065
$ws.Cells.Item(11, 40).Value = 0.255
$ws.Cells.Item(11, 41).Value = 0.891

# Row 12
$ws.Cells.Item(12, 2).Value = 1.389
$ws.Cells.Item(12, 3).Value = 0.682
$ws.Cells.Item(12, 4).Value = 0.826
$ws.Cells.Item(12, 5).Value = 1.633
$ws.Cells.Item(12, 6).Value = 1.032
$ws.Cells.Item(12, 7).Value = 1.016
$ws.Cells.Item(12, 8).Value = 1.556
$ws.Cells.Item(12, 9).Value = 1.191
$ws.Cells.Item(12, 10).Value = 1.091
$ws.Cells.Item(12, 11).Value = 1.407
$ws.Cells.Item(12, 12).Value = 0.538
$ws.Cells.Item(12, 13).Value = 0.733
$ws.Cells.Item(12, 14).Value = 1.343
$ws.Cells.Item(12, 15).Value = 0.511
$ws.Cells.Item(12, 16).Value = 0.715
$ws.Cells.Item(12, 26).Value = 1.25
$ws.Cells.Item(12, 27).Value = 0.299
$ws.Cells.Item(12, 28).Value = 0.546
$ws.Cells.Item(12, 29).Value = 2.032
$ws.Cells.Item(12, 30).Value = 3.902
$ws.Cells.Item(12, 31).Value = 1.975
$ws.Cells.Item(12, 32).Value = 1.22
$ws.Cells.Item(12, 33).Value = 0.22
$ws.Cells.Item(12, 34).Value = 0.469
$ws.Cells.Item(12, 35).Value = 1.029
$ws.Cells.Item(12, 37).Value = 0.169
$ws.Cells.Item(12, 38).Value = 1.1
$ws.Cells.Item(12, 39).Value = 0.09
$ws.Cells.Item(12, 40).Value = 0.3
$ws.Cells.Item(12, 41).Value = 1.116

# Row 13
$ws.Cells.Item(13, 2).Value = 3.465
$ws.Cells.Item(13, 3).Value = 1.365
$ws.Cells.Item(13, 4).Value = 1.168
$ws.Cells.Item(13, 5).Value = 4.541
$ws.Cells.Item(13, 6).Value = 0.735
$ws.Cells.Item(13, 7).Value = 0.857
$ws.Cells.Item(13, 8).Value = 4.5
$ws.Cells.Item(13, 9).Value = 0.95
$ws.Cells.Item(13, 10).Value = 0.975
$ws.Cells.Item(13, 11).Value = 2.3
$ws.Cells.Item(13, 12).Value = 0.61
$ws.Cells.Item(13, 13).Value = 0.781
$ws.Cells.Item(13, 14).Value = 3.302
$ws.Cells.Item(13, 15).Value = 0.769
$ws.Cells.Item(13, 16).Value = 0.877
$ws.Cells.Item(13, 26).Value = 2.833
$ws.Cells.Item(13, 27).Value = 3.901
$ws.Cells.Item(13, 28).Value = 1.975
$ws.Cells.Item(13, 29).Value = 6.286
$ws.Cells.Item(13, 30).Value = 2.966
$ws.Cells.Item(13, 31).Value = 1.722
$ws.Cells.Item(13, 32).Value = 1.628
$ws.Cells.Item(13, 33).Value = 0.699
$ws.Cells.Item(13, 34).Value = 0.836
$ws.Cells.Item(13, 35).Value = 1.233
$ws.Cells.Item(13, 36).Value = 0.178
$ws.Cells.Item(13, 37).Value = 0.422
$ws.Cells.Item(13, 38).Value = 1.651
$ws.Cells.Item(13, 39).Value = 0.785
$ws.Cells.Item(13, 40).Value = 0.886
$ws.Cells.Item(13, 41).Value = 1.504
